$d = $word.ActiveDocument

function Get-ParagraphByStyle([string]$styleName) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Style.NameLocal -eq $styleName) {
            return $p
        }
    }
    return $null
}

function Set-ParagraphRunXml($paragraph, [string]$text) {
    # Build a Document-level Range (not a Paragraph.Range) spanning the
    # paragraph's text (excluding the trailing paragraph mark) so that
    # InsertXML fully replaces the existing runs (instead of appending
    # after them), collapsing the paragraph down to a single run while
    # keeping the paragraph mark / pPr (and its style) untouched.
    $pr = $paragraph.Range
    $start = $pr.Start
    $end = $pr.End - 1
    $r = $d.Range($start, $end)

    $escaped = $text.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p><w:r><w:t xml:space="preserve">' + $escaped + '</w:t></w:r></w:p></w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml)
}

# Title paragraph: "First" " " "VIP" " " "student" " " "team" " " "outputs"
#   -> single run "First VIP student team outputs"
Set-ParagraphRunXml (Get-ParagraphByStyle "Title") "First VIP student team outputs"

# Author paragraph: "Tom" " " "Coleman" -> single run "Tom Coleman"
Set-ParagraphRunXml (Get-ParagraphByStyle "Author") "Tom Coleman"

# Abstract paragraph -> single merged run
Set-ParagraphRunXml (Get-ParagraphByStyle "Abstract") "Materials from our first VIP student teams are available now, with many more guides to follow!"
